{"js": "const paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nconst p = paras.items[5];\np.font.color = \"Accent 6\";\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$p = $d.Paragraphs(6)\n$r = $p.Range\n$r.Font.TextColor.ObjectThemeColor = 6\n"}
